# Section 13: Animations completed
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Unity Course Video Listing")

# Mark all the "Animation" section videos (rows 174-192) as completed by
# filling in their COMPLETED date (column G) with the same date already
# present in the TARGET DATE column (H) for each of those rows: 3/2/2017
# (Excel serial date 42796).
for ($r = 174; $r -le 192; $r++) {
    $ws.Cells.Item($r, 7).Value = 42796
}

# Add completion notes for two of the rows.
$ws.Range("K188").Value = "animations weren't working since he used a different avatar from the original"
$ws.Range("K189").Value = "watched videos and took quizzes but didn" + [char]0x2019 + "t follow along in unity"

# Update the active selection to reflect where the author ended up after
# making the edits.
$ws.Range("K190").Select()
